# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.323.47'
$ws.Range("E2").Value = '  -1.98%  '

# Row 3
$ws.Range("D3").Value = '2.614.96'
$ws.Range("E3").Value = '  -3.96%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '550.78'
$ws.Range("E5").Value = '  -1.79%  '

# Row 6
$ws.Range("D6").Value = '153.88'
$ws.Range("E6").Value = '  -3.58%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.590'
$ws.Range("E8").Value = '  -0.30%  '

# Row 9
$ws.Range("E9").Value = '  -2.89%  '

# Row 10
$ws.Range("E10").Value = '  -3.99%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.48'
$ws.Range("E11").Value = '  -0.99%  '

# Row 12
$ws.Range("E12").Value = '  -2.16%  '

# Row 13
$ws.Range("D13").Value = '3.080.87'
$ws.Range("E13").Value = '  -3.83%  '

# Row 14
$ws.Range("D14").Value = '25.68'
$ws.Range("E14").Value = '  -3.09%  '

# Row 15
$ws.Range("D15").Value = '62.235.37'

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000144'
$ws.Range("E16").Value = '  -2.53%  '

# Row 17
$ws.Range("D17").Value = '2.618.16'
$ws.Range("E17").Value = '  -3.56%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.61'
$ws.Range("E18").Value = '  -4.98%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.53'
$ws.Range("E19").Value = '  -2.84%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '340.33'
$ws.Range("E20").Value = '  -3.23%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.10'
$ws.Range("E21").Value = '  -5.91%  '

# Row 22
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.06%  '

# Row 23
$ws.Range("E23").Value = '  -2.99%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.77'
$ws.Range("E24").Value = '  -1.55%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.168'
$ws.Range("E25").Value = '  -0.11%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.03%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.03'
$ws.Range("E27").Value = '  -1.95%  '

# Row 28
$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.35'
$ws.Range("E28").Value = '  +0.33%  '

# Row 29
$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.19'
$ws.Range("E29").Value = '  +1.28%  '

# Row 30
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0825'
$ws.Range("E30").Value = '  -7.39%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.90'
$ws.Range("E31").Value = '  -2.93%  '

# Row 32
$ws.Range("E32").Value = '  +0.01%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '159.58'
$ws.Range("E33").Value = '  -4.08%  '

# Row 34
$ws.Range("D34").Value = '19.23'
$ws.Range("E34").Value = '  -3.11%  '

# Row 35
$ws.Range("E35").Value = '  -3.29%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.42'
$ws.Range("E36").Value = '  -3.64%  '

# Row 37
$ws.Range("D37").Value = '1.74'
$ws.Range("E37").Value = '  -2.50%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '337.98'
$ws.Range("E38").Value = '  -1.23%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.10'
$ws.Range("E39").Value = '  -1.91%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.889'
$ws.Range("E40").Value = '  -6.36%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '37.73'
$ws.Range("E41").Value = '  -1.70%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.89'
$ws.Range("E42").Value = '  -3.64%  '

# Row 43
$ws.Range("E43").Value = '  +0.18%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.609'
$ws.Range("E44").Value = '  -2.25%  '

# Row 45
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.37'
$ws.Range("E45").Value = '  -4.48%  '

# Row 46
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.131.44'
$ws.Range("E46").Value = '  +1.88%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.97'
$ws.Range("E47").Value = '  -0.91%  '

# Row 48
$ws.Range("E48").Value = '  -4.66%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0547'
$ws.Range("E49").Value = '  -4.79%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0962'
$ws.Range("E50").Value = '  -2.54%  '

# Row 51
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0237'
$ws.Range("E51").Value = '  -3.61%  '
